# Exchange rate.xlsx - CreateMyCells finally work. Change headers cells style bold=true.
#
# 1. The header row (row 1: A1:D1) gets a bold font style.
# 2. The combined shared string "[no, effectiveDate, mid]" that was (incorrectly)
#    repeated across D2, E2 and F2 is replaced by the three separate column
#    header labels "no", "effectiveDate" and "mid".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the bogus combined header string into its three real values.
$ws.Range("D2").Value = "no"
$ws.Range("E2").Value = "effectiveDate"
$ws.Range("F2").Value = "mid"

# Make the header cells (table/currency/code/rates) bold.
$ws.Range("A1:D1").Font.Bold = $true
